$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: content-only fix (Objetivos: body text was in the wrong place) ---
$ws.Range("B10").Value = 'Apresentar o formalismo para descrição de sistemas quânticos. Estudar diversas aplicações da equação de Schroedinger independente do tempo. Descrever a estrutura eletrônica de átomos e moléculas.'
$ws.Range("C10").Value = 'Apresentar o formalismo para descrição de sistemas quânticos. Estudar diversas aplicações da equação de Schroedinger independente do tempo. Descrever a estrutura eletrônica de átomos e moléculas.'

# --- Rebuild rows 13-25 into the new 13-28 layout ---
# The old rows 13:25 interleave data incorrectly (extra rows, shifted columns).
# Remove them and open up 16 fresh rows at 13 so the sheet grows to 28 rows,
# matching the corrected layout (dimension A1:C28).
$ws.Rows("13:25").Delete()
$ws.Rows("13:28").Insert()

# Row-insert carries column-A formatting down from row 12; strip it back
# off the rows that should have no entry in column A.
foreach ($addr in @("A13","A14","A15","A26","A27","A28")) {
    $ws.Range($addr).Style = "Normal"
}

# Apply the bold "label" style (col A), normal wrap style (col B) and
# red wrap style (col C) used throughout the sheet, by copying formats
# from existing template cells of each kind.
$ws.Range("A12").Copy() | Out-Null
foreach ($addr in @("A16","A17","A18","A19","A20","A21","A22","A23","A24","A25")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$ws.Range("B10").Copy() | Out-Null
foreach ($addr in @("B13","B14","B15","B16","B17","B18","B19","B21","B22","B23","B24","B26","B27","B28")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$ws.Range("C10").Copy() | Out-Null
foreach ($addr in @("C13","C14","C15","C16","C17","C18","C19","C21","C22","C23","C24","C26","C27","C28")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Row 13
$ws.Range("B13").Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Range("C13").Value = '5840730 - Antonio Jefferson da Silva Machado'

# Row 14
$ws.Range("B14").Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range("C14").Value = '1176388 - Luiz Tadeu Fernandes Eleno'

# Row 15
$ws.Range("B15").Value = '1341653 - Maria José Ramos Sandim'
$ws.Range("C15").Value = '1341653 - Maria José Ramos Sandim'

# Row 16
$ws.Range("A16").Value = 'Programa resumido:'
$ws.Range("B16").Value = 'Introdução aos conceitos da Mecânica Quântica. • Ferramentas matemáticas da Mecânica Quântica. A equação de Schroedinger e aplicações unidimensionais e tridimensionais. Problemas em coordenadas retangulares. Problemas em coordenadas esféricas. Átomos com um elétron. Teoria geral. • Propriedades gerais do momento angular.'
$ws.Range("C16").Value = 'Introdução aos conceitos da Mecânica Quântica. • Ferramentas matemáticas da Mecânica Quântica. A equação de Schroedinger e aplicações unidimensionais e tridimensionais. Problemas em coordenadas retangulares. Problemas em coordenadas esféricas. Átomos com um elétron. Teoria geral. • Propriedades gerais do momento angular.'
$ws.Rows(16).RowHeight = 60

# Row 17
$ws.Range("A17").Value = 'Short syllabus:'
$ws.Range("B17").Value = '• Introduction to the concepts of Quantum Mechanics. • Mathematical tools of Quantum Mechanics. • The Schrödinger equation and one- and three-dimensional applications. • Quantum formalism. • Problems in rectangular coordinates and spherical coordinates. • Hydrogen atoms and orbitals. • General properties of angular momentum. • Spin. • Fermions and bosons.'
$ws.Range("C17").Value = '• Introduction to the concepts of Quantum Mechanics. • Mathematical tools of Quantum Mechanics. • The Schrödinger equation and one- and three-dimensional applications. • Quantum formalism. • Problems in rectangular coordinates and spherical coordinates. • Hydrogen atoms and orbitals. • General properties of angular momentum. • Spin. • Fermions and bosons.'
$ws.Rows(17).RowHeight = 60

# Row 18
$ws.Range("A18").Value = 'Programa:'
$ws.Range("B18").Value = '• Equação de Schrödinger. • Função de onda e interpretação estatística da mecânica quântica. • Valores esperados e operadores. Os operadores posição e momento; operadores energia cinética e potencial; o operador Hamiltoniano. • A equação de Schrödinger independente do tempo. Separação de variáveis e estados estacionários. • Aplicações unidimensionais:  poço quadrado infinito; oscilador harmônico; partícula livre;  transformada de Fourier e sua relação com o princípio da incerteza de Heisenberg; Poços e barreiras de potencial. • Formalismo quântico: opserváveis e operadores hermitianos. Estados determinados, autoestados e autovalores de operadores hermitianos. Base de autoestados; interpretação estatística generalizada: medidas de observáveis e suas probabilidades. Comutadores e operadores que compartilham autoestados ; princípio da incerteza generalizado. • Mecânica Quântica em três dimensões. • Átomo de hidrogênio: modelo de Bohr e o número quântico principal. Solução completa e os demais números quânticos. • Coordenadas esféricas e Momento angular.  • Momento angulas de spin. • Problemas de muitos corpos. • Partículas idênticas: férmions e bósons.'
$ws.Range("C18").Value = '• Equação de Schrödinger. • Função de onda e interpretação estatística da mecânica quântica. • Valores esperados e operadores. Os operadores posição e momento; operadores energia cinética e potencial; o operador Hamiltoniano. • A equação de Schrödinger independente do tempo. Separação de variáveis e estados estacionários. • Aplicações unidimensionais:  poço quadrado infinito; oscilador harmônico; partícula livre;  transformada de Fourier e sua relação com o princípio da incerteza de Heisenberg; Poços e barreiras de potencial. • Formalismo quântico: opserváveis e operadores hermitianos. Estados determinados, autoestados e autovalores de operadores hermitianos. Base de autoestados; interpretação estatística generalizada: medidas de observáveis e suas probabilidades. Comutadores e operadores que compartilham autoestados ; princípio da incerteza generalizado. • Mecânica Quântica em três dimensões. • Átomo de hidrogênio: modelo de Bohr e o número quântico principal. Solução completa e os demais números quânticos. • Coordenadas esféricas e Momento angular.  • Momento angulas de spin. • Problemas de muitos corpos. • Partículas idênticas: férmions e bósons.'
$ws.Rows(18).RowHeight = 120

# Row 19
$ws.Range("A19").Value = 'Syllabus:'
$ws.Range("B19").Value = '• Schrödinger''s equation. • Wave function and statistical interpretation of quantum mechanics. • Expected values and operators. The position and moment operators; kinetic and potential energy operators; the Hamiltonian operator. • The time-independent Schrödinger equation. Separation of variables and steady states. • One-dimensional applications: infinite square well; harmonic oscillator; free particle; Fourier transform and its relationship with the Heisenberg uncertainty principle; Potential square wells and barriers. • Quantum formalism: hermitian operators and observables. Determined states, eigenstates and eigenvalues of Hermitian operators. Basis of Eigenstates; generalized statistical interpretation: measures of observables and their probabilities. Comutators and operators that share eigenstates; generalized uncertainty principle. • Quantum Mechanics in three dimensions. • Hydrogen atom: Bohr model and the principal quantum number. Complete solution and the other quantum numbers. • Spherical coordinates and Angular momentum. • Spin angular momentum. • Many-body problems. • Identical particles: fermions and bosons.'
$ws.Range("C19").Value = '• Schrödinger''s equation. • Wave function and statistical interpretation of quantum mechanics. • Expected values and operators. The position and moment operators; kinetic and potential energy operators; the Hamiltonian operator. • The time-independent Schrödinger equation. Separation of variables and steady states. • One-dimensional applications: infinite square well; harmonic oscillator; free particle; Fourier transform and its relationship with the Heisenberg uncertainty principle; Potential square wells and barriers. • Quantum formalism: hermitian operators and observables. Determined states, eigenstates and eigenvalues of Hermitian operators. Basis of Eigenstates; generalized statistical interpretation: measures of observables and their probabilities. Comutators and operators that share eigenstates; generalized uncertainty principle. • Quantum Mechanics in three dimensions. • Hydrogen atom: Bohr model and the principal quantum number. Complete solution and the other quantum numbers. • Spherical coordinates and Angular momentum. • Spin angular momentum. • Many-body problems. • Identical particles: fermions and bosons.'
$ws.Rows(19).RowHeight = 120

# Row 20
$ws.Range("A20").Value = 'Avaliação:'

# Row 21
$ws.Range("A21").Value = 'Método:'
$ws.Range("B21").Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Range("C21").Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Rows(21).RowHeight = 60

# Row 22
$ws.Range("A22").Value = 'Critério:'
$ws.Range("B22").Value = 'Média aritmética de três provas: P1 (peso 1), P2 (peso 1) e P3 (peso 2).'
$ws.Range("C22").Value = 'Média aritmética de três provas: P1 (peso 1), P2 (peso 1) e P3 (peso 2).'
$ws.Rows(22).RowHeight = 60

# Row 23
$ws.Range("A23").Value = 'Norma de recuperação:'
$ws.Range("B23").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C23").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Rows(23).RowHeight = 60

# Row 24
$ws.Range("A24").Value = 'Bibliografia:'
$ws.Range("B24").Value = 'Griffiths, D. J. Mecânica Quântica, 2a ed., Pearson, 2011.ZETTILI, N. Quantum Mechanics: Concepts and Applications, Wiley, 2009.CLAUDE COHEN-TANNOUDJI, BERNARD DIU, FRANK LALOE. Quantum Mechanics, Vol 1 e 2. Ed. John Wiley and Sons, 1987.GASIOROWICZ, S., Física Quântica, Guanabara Dois, RJ. 1979.FEYNMAN, R.P., LEIGHTON, R.B. AND SANDS, M., The Feynman Lectures on Physics, vol.3, Addison-Wesley, 1975.MERZBACHER, E., Quantum Mechanics, John Wiley & Sons, Nova Iorque, 1970.EISBERG, R.; RESNICK, R., Física Quântica, Átomos, Moléculas, Sólidos, Núcleos e Partículas, Ed. Campus, 1978.'
$ws.Range("C24").Value = 'Griffiths, D. J. Mecânica Quântica, 2a ed., Pearson, 2011.ZETTILI, N. Quantum Mechanics: Concepts and Applications, Wiley, 2009.CLAUDE COHEN-TANNOUDJI, BERNARD DIU, FRANK LALOE. Quantum Mechanics, Vol 1 e 2. Ed. John Wiley and Sons, 1987.GASIOROWICZ, S., Física Quântica, Guanabara Dois, RJ. 1979.FEYNMAN, R.P., LEIGHTON, R.B. AND SANDS, M., The Feynman Lectures on Physics, vol.3, Addison-Wesley, 1975.MERZBACHER, E., Quantum Mechanics, John Wiley & Sons, Nova Iorque, 1970.EISBERG, R.; RESNICK, R., Física Quântica, Átomos, Moléculas, Sólidos, Núcleos e Partículas, Ed. Campus, 1978.'
$ws.Rows(24).RowHeight = 120

# Row 25
$ws.Range("A25").Value = 'Requisitos:'

# Row 26
$ws.Range("B26").Value = 'LOB1021 -  Física IV  (Requisito)
'
$ws.Range("C26").Value = 'LOB1021 -  Física IV  (Requisito)
'
$ws.Rows(26).RowHeight = 30

# Row 27
$ws.Range("B27").Value = 'LOM3253 -  Física Matemática  (Requisito)
'
$ws.Range("C27").Value = 'LOM3253 -  Física Matemática  (Requisito)
'
$ws.Rows(27).RowHeight = 30

# Row 28
$ws.Range("B28").Value = 'LOM3257 -  Mecânica Clássica  (Requisito fraco)
'
$ws.Range("C28").Value = 'LOM3257 -  Mecânica Clássica  (Requisito fraco)
'
$ws.Rows(28).RowHeight = 30

Write-Output "done"
